$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pythonCode")

# The cell A1 on the "pythonCode" sheet contained "pCode"; rename it to "pyCode".
$ws.Range("A1").Value = "pyCode"

# Activate this sheet and scroll/select back to the top-left (A1) so the
# saved view no longer points at A10 with that cell selected.
$ws.Activate()
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
